$d = $word.ActiveDocument

$targets = @(
    "This is an annotatable resource in the casebook.",
    "highlighted: content to highlight; elided: content to elide; replaced: content to replace; linked: content to link; noted: content to note; highlighted2: second highlight content;",
    "This is the second chapter of the casebook."
)

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text
    # Paragraph.Range.Text includes the trailing paragraph mark (CR, or CR+LF
    # when the run itself also embeds a literal line-feed character).
    $trimmed = $text.TrimEnd([char]13, [char]10)
    # Collapse any embedded literal line-feed characters to single spaces,
    # joining wrapped lines back into one.
    $collapsed = $trimmed.Replace([char]10, ' ')

    foreach ($target in $targets) {
        if ($collapsed -eq $target) {
            $p.Range.Text = $target
            break
        }
    }
}
